$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New label / value cells (rows 6, 8, 9, 10)
# ------------------------------------------------------------------
$ws.Range("B6").Value = "Descripción del Activo"

$ws.Range("B8").Value = "Valor del activo"
$ws.Range("C8").Value = 56000

$ws.Range("B9").Value = "Vida util"
$ws.Range("C9").Value = 5

$ws.Range("B10").Value = "Descripción"
$ws.Range("C10").Value = "PC gamer ultima generación"

# Bold, yellow-filled labels (B6, B8, B9, B10)
foreach ($addr in @("B6","B8","B9","B10")) {
    $ws.Range($addr).Font.Bold = $true
}
foreach ($addr in @("B6","B8","B9","B10")) {
    $ws.Range($addr).Interior.Color = 65535
}
foreach ($addr in @("B6","B8","B9","B10")) {
    $ws.Range($addr).Interior.PatternColor = 0
}

# Centered values next to the bold labels (C8, C9, C10)
foreach ($addr in @("C8","C9","C10")) {
    $ws.Range($addr).HorizontalAlignment = -4108
}
foreach ($addr in @("C8","C9","C10")) {
    $ws.Range($addr).VerticalAlignment = -4107
}

$ws.Rows(6).RowHeight = 15

# ------------------------------------------------------------------
# Update the sample employee/asset row (row 3)
# ------------------------------------------------------------------
$ws.Range("B3").Value = "PC"
$ws.Range("C3").Value = "Kevin"
$ws.Range("D3").Value = "SflpybZh"
$ws.Range("E3").Value = "'002"
$ws.Range("E3").ClearFormats()

# ------------------------------------------------------------------
# Re-style the header row (row 2): bold, size 16, centered, gray0625
# pattern over a yellow/black fill
# ------------------------------------------------------------------
$ws.Rows(2).ClearFormats()

$hdr = $ws.Range("B2:E2")
$hdr.Font.Bold = $true
$hdr.Font.Size = 16
$hdr.Interior.Color = 65535
$hdr.Interior.Pattern = 18
$hdr.Interior.PatternColor = 0
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4107
